$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 2 (pushes the existing 20 data
# rows from rows 2-21 down to rows 5-24).
$ws.Rows("2:4").Insert()

# Excel's Insert() carries formatting down from the row above into the new
# rows; the source data has no special formatting on these cells, so clear
# it back to the default style.
$ws.Range("A2:C4").ClearFormats()

# Fill the 3 newly inserted rows with their values.
$ws.Cells.Item(2, 1).Value = -1.012200951576233
$ws.Cells.Item(2, 2).Value = -1.622797250747681
$ws.Cells.Item(2, 3).Value = -0.285490870475769

$ws.Cells.Item(3, 1).Value = 0.4157860279083252
$ws.Cells.Item(3, 2).Value = 2.06141996383667
$ws.Cells.Item(3, 3).Value = -0.3307646214962005

$ws.Cells.Item(4, 1).Value = 1.623263239860535
$ws.Cells.Item(4, 2).Value = 1.307212591171265
$ws.Cells.Item(4, 3).Value = -1.966478228569031

# Append 7 new rows after the (now shifted) last data row (row 24), i.e.
# rows 25-31.
$ws.Cells.Item(25, 1).Value = 3.35325288772583
$ws.Cells.Item(25, 2).Value = 15.64194393157959
$ws.Cells.Item(25, 3).Value = -7.271495342254639

$ws.Cells.Item(26, 1).Value = 7.133343696594238
$ws.Cells.Item(26, 2).Value = -10.52308177947998
$ws.Cells.Item(26, 3).Value = -2.896986722946167

$ws.Cells.Item(27, 1).Value = 0.8908939957618713
$ws.Cells.Item(27, 2).Value = -0.6539392471313477
$ws.Cells.Item(27, 3).Value = 2.433596611022949

$ws.Cells.Item(28, 1).Value = -1.749364018440247
$ws.Cells.Item(28, 2).Value = -0.2597913742065429
$ws.Cells.Item(28, 3).Value = 1.243696212768555

$ws.Cells.Item(29, 1).Value = -14.04837512969971
$ws.Cells.Item(29, 2).Value = 1.555951952934265
$ws.Cells.Item(29, 3).Value = 13.38611316680908

$ws.Cells.Item(30, 1).Value = 8.504871368408203
$ws.Cells.Item(30, 2).Value = -7.031145095825195
$ws.Cells.Item(30, 3).Value = 0.2657834887504577

$ws.Cells.Item(31, 1).Value = 3.329284429550171
$ws.Cells.Item(31, 2).Value = -2.393516063690185
$ws.Cells.Item(31, 3).Value = 2.149171113967896
